# Updates market-price derived cells across all Leve-profit worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market data
# pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 364.85715
$ws.Range("I33").Value = 355.85715
$ws.Range("K33").Value = 355.85715
$ws.Range("M33").Value = -126.85715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 6628.913
$ws.Range("I98").Value = 6565.268
$ws.Range("K98").Value = 6565.268
$ws.Range("M98").Value = -5067.268

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 6628.913
$ws.Range("I122").Value = 6565.268
$ws.Range("K122").Value = 19695.804
$ws.Range("M122").Value = -17245.804

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1988
$ws.Range("I132").Value = 1499.7333
$ws.Range("J132").Value = 3819
$ws.Range("K132").Value = 4499.199900000001
$ws.Range("L132").Value = 11457
$ws.Range("M132").Value = -1969.199900000001
$ws.Range("N132").Value = -16517

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 108460.664
$ws.Range("J133").Value = 108460.664
$ws.Range("L133").Value = 108460.664
$ws.Range("N133").Value = -118580.664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 905.6667
$ws.Range("I135").Value = 905.6667
$ws.Range("K135").Value = 8151.0003
$ws.Range("M135").Value = -5616.0003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1281246.8
$ws.Range("I138").Value = 6578.6665
$ws.Range("J138").Value = 1686822.9
$ws.Range("K138").Value = 19735.9995
$ws.Range("L138").Value = 5060468.699999999
$ws.Range("M138").Value = -14595.9995
$ws.Range("N138").Value = -5070748.699999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6766.593
$ws.Range("I141").Value = 4459.6665
$ws.Range("K141").Value = 13378.9995
$ws.Range("M141").Value = -8198.999500000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 999
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 998
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 998
$ws.Range("M17").Value = -827
$ws.Range("N17").Value = -1344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 113386.25
$ws.Range("I32").Value = 158020.1
$ws.Range("J32").Value = 15191.8
$ws.Range("K32").Value = 158020.1
$ws.Range("L32").Value = 15191.8
$ws.Range("M32").Value = -157733.1
$ws.Range("N32").Value = -15765.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5745.769
$ws.Range("I61").Value = 4924.5264
$ws.Range("K61").Value = 4924.5264
$ws.Range("M61").Value = -4712.5264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2713.2593
$ws.Range("I74").Value = 2356.8462
$ws.Range("K74").Value = 2356.8462
$ws.Range("M74").Value = -1482.8462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2713.2593
$ws.Range("I77").Value = 2356.8462
$ws.Range("K77").Value = 11784.231
$ws.Range("M77").Value = -7416.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 945.73334
$ws.Range("I110").Value = 881.2692
$ws.Range("J110").Value = 1364.75
$ws.Range("K110").Value = 881.2692
$ws.Range("L110").Value = 1364.75
$ws.Range("M110").Value = 1163.7308
$ws.Range("N110").Value = -5454.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3106.0588
$ws.Range("I132").Value = 2117.0715
$ws.Range("K132").Value = 6351.2145
$ws.Range("M132").Value = -3821.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5745.769
$ws.Range("I136").Value = 4924.5264
$ws.Range("K136").Value = 14773.5792
$ws.Range("M136").Value = -12223.5792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 12408.357
$ws.Range("J105").Value = 5814.375
$ws.Range("L105").Value = 5814.375
$ws.Range("N105").Value = -9308.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3219.5098
$ws.Range("I107").Value = 2237.375
$ws.Range("J107").Value = 6790.909
$ws.Range("K107").Value = 2237.375
$ws.Range("L107").Value = 6790.909
$ws.Range("M107").Value = -317.375
$ws.Range("N107").Value = -10630.909

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2513.0833
$ws.Range("I134").Value = 2528.05
$ws.Range("K134").Value = 7584.150000000001
$ws.Range("M134").Value = -5049.150000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4058.1428
$ws.Range("I31").Value = 1675.6923
$ws.Range("K31").Value = 1675.6923
$ws.Range("M31").Value = -1380.6923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4058.1428
$ws.Range("I34").Value = 1675.6923
$ws.Range("K34").Value = 1675.6923
$ws.Range("M34").Value = -1473.6923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2340.0908
$ws.Range("I58").Value = 2380.3
$ws.Range("J58").Value = 1938
$ws.Range("K58").Value = 2380.3
$ws.Range("L58").Value = 1938
$ws.Range("M58").Value = -2177.3
$ws.Range("N58").Value = -2344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2380.3
$ws.Range("J132").Value = 2665.4285
$ws.Range("L132").Value = 7996.2855
$ws.Range("N132").Value = -13056.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 49329.285
$ws.Range("I134").Value = 56134.168
$ws.Range("K134").Value = 168402.504
$ws.Range("M134").Value = -165867.504

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2340.0908
$ws.Range("I136").Value = 2380.3
$ws.Range("J136").Value = 1938
$ws.Range("K136").Value = 7140.900000000001
$ws.Range("L136").Value = 5814
$ws.Range("M136").Value = -4590.900000000001
$ws.Range("N136").Value = -10914

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 225405.08
$ws.Range("J141").Value = 265486.1
$ws.Range("L141").Value = 265486.1
$ws.Range("N141").Value = -275846.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3760175.8
$ws.Range("I4").Value = 1635638.2
$ws.Range("K4").Value = 4906914.6
$ws.Range("M4").Value = -4906802.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2193.75
$ws.Range("I63").Value = 925
$ws.Range("K63").Value = 2775
$ws.Range("M63").Value = -2026

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 2193.75
$ws.Range("I66").Value = 925
$ws.Range("K66").Value = 8325
$ws.Range("M66").Value = -4581

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 303811.8
$ws.Range("I121").Value = 1030
$ws.Range("K121").Value = 3090
$ws.Range("M121").Value = -1780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2334.1365
$ws.Range("I140").Value = 2434.4375
$ws.Range("K140").Value = 7303.3125
$ws.Range("M140").Value = -2123.3125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22226516
$ws.Range("I70").Value = 47621960
$ws.Range("K70").Value = 47621960
$ws.Range("M70").Value = -47621690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 22226516
$ws.Range("I73").Value = 47621960
$ws.Range("K73").Value = 47621960
$ws.Range("M73").Value = -47621024

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 905.0454999999999
$ws.Range("I97").Value = 816.7895
$ws.Range("K97").Value = 816.7895
$ws.Range("M97").Value = -320.7895

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 20835308
$ws.Range("I102").Value = 26316874
$ws.Range("K102").Value = 26316874
$ws.Range("M102").Value = -26315252

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7837
$ws.Range("I122").Value = 6898.5
$ws.Range("K122").Value = 20695.5
$ws.Range("M122").Value = -18245.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3781.125
$ws.Range("I126").Value = 2518.75
$ws.Range("K126").Value = 7556.25
$ws.Range("M126").Value = -5086.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2202.125
$ws.Range("I132").Value = 1833.2222
$ws.Range("J132").Value = 2503.9546
$ws.Range("K132").Value = 5499.6666
$ws.Range("L132").Value = 7511.8638
$ws.Range("M132").Value = -2969.6666
$ws.Range("N132").Value = -12571.8638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3209.2
$ws.Range("I7").Value = 2298.8572
$ws.Range("K7").Value = 2298.8572
$ws.Range("M7").Value = -2186.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1401.8889
$ws.Range("I16").Value = 1462.6471
$ws.Range("J16").Value = 369
$ws.Range("K16").Value = 1462.6471
$ws.Range("L16").Value = 369
$ws.Range("M16").Value = -1292.6471
$ws.Range("N16").Value = -709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 70000
$ws.Range("J36").Value = 70000
$ws.Range("L36").Value = 70000
$ws.Range("N36").Value = -71124

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6899.1333
$ws.Range("I100").Value = 3370.1428
$ws.Range("K100").Value = 3370.1428
$ws.Range("M100").Value = -2829.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3209.2
$ws.Range("I126").Value = 2298.8572
$ws.Range("K126").Value = 6896.571599999999
$ws.Range("M126").Value = -4426.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3626.318
$ws.Range("I132").Value = 3385.6667
$ws.Range("K132").Value = 10157.0001
$ws.Range("M132").Value = -7627.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 33219.35
$ws.Range("I136").Value = 2944.182
$ws.Range("J136").Value = 70222.336
$ws.Range("K136").Value = 8832.545999999998
$ws.Range("L136").Value = 210667.008
$ws.Range("M136").Value = -6282.545999999998
$ws.Range("N136").Value = -215767.008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 91836.664
$ws.Range("J94").Value = 91836.664
$ws.Range("L94").Value = 91836.664
$ws.Range("N94").Value = -93638.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 983.44446
$ws.Range("I113").Value = 869.4706
$ws.Range("J113").Value = 1177.2
$ws.Range("K113").Value = 2608.4118
$ws.Range("L113").Value = 3531.6
$ws.Range("M113").Value = -438.4117999999999
$ws.Range("N113").Value = -7871.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3477.9333
$ws.Range("I132").Value = 2357.2856
$ws.Range("J132").Value = 7400.2
$ws.Range("K132").Value = 7071.8568
$ws.Range("L132").Value = 22200.6
$ws.Range("M132").Value = -4541.8568
$ws.Range("N132").Value = -27260.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 70487.336
$ws.Range("J133").Value = 70487.336
$ws.Range("L133").Value = 70487.336
$ws.Range("N133").Value = -80607.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3643.2856
$ws.Range("I136").Value = 1062.4546
$ws.Range("K136").Value = 3187.3638
$ws.Range("M136").Value = -637.3638000000001
